$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Force the "Price" cells to remain plain Text (matching the source data,
# which stores prices like "61.058.58" / "10.60" / "581.00" as text, not numbers)
$dPriceCells = @("D2","D3","D4","D5","D6","D9","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D23","D26","D28","D31","D34","D38","D39","D40","D43","D46","D47","D50")
foreach ($addr in $dPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.058.58"
$ws.Range("D3").Value = "2.434.21"
$ws.Range("D4").Value = "0.998"
$ws.Range("D5").Value = "572.29"
$ws.Range("D6").Value = "140.64"
$ws.Range("D9").Value = "2.421.74"
$ws.Range("D12").Value = "5.14"
$ws.Range("D13").Value = "0.340"
$ws.Range("D14").Value = "26.10"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("D16").Value = "2.856.35"
$ws.Range("D17").Value = "61.040.68"
$ws.Range("D18").Value = "2.440.13"
$ws.Range("D19").Value = "10.60"
$ws.Range("D20").Value = "7.29"
$ws.Range("D21").Value = "324.19"
$ws.Range("D23").Value = "6.16"
$ws.Range("D26").Value = "65.24"
$ws.Range("D28").Value = "581.00"
$ws.Range("D31").Value = "7.89"
$ws.Range("D34").Value = "0.132"
$ws.Range("D38").Value = "150.72"
$ws.Range("D39").Value = "1.39"
$ws.Range("D40").Value = "18.26"
$ws.Range("D43").Value = "41.70"
$ws.Range("D46").Value = "0.0₆0285"
$ws.Range("D47").Value = "141.73"
$ws.Range("D50").Value = "19.63"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  -2.01%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -3.79%  "
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("E27").Value = "  -5.40%  "
$ws.Range("E28").Value = "  -6.32%  "
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  -4.24%  "
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("E32").Value = "  -5.43%  "
$ws.Range("E33").Value = "  -2.13%  "
$ws.Range("E34").Value = "  -6.12%  "
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("E36").Value = "  -6.11%  "
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("E39").Value = "  -3.44%  "
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("E44").Value = "  -5.89%  "
$ws.Range("E45").Value = "  -4.77%  "
$ws.Range("E46").Value = "  +25.28%  "
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("E48").Value = "  -2.83%  "
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("E51").Value = "  -3.37%  "
